# Update datatype column (H) values from "date" to "datetime"
# for the rows that describe date-time fields.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$rows = @(6, 12, 45, 49, 52, 73, 76, 80, 88, 90)

foreach ($r in $rows) {
    $cell = $ws.Range("H$r")
    if ($cell.Text -eq "date") {
        $cell.Value = "datetime"
    }
}
